# Department meetings template:
# 1. Show the submitter's full name instead of login name everywhere the
#    sales person is referenced on the generated report.
# 2. (Region/department required-input validation is enforced elsewhere in
#    the app; nothing in this template's cell layout needs to change for it.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "D2" placeholder cell drives the column that used to resolve to the
# user's login name - point it at the full-name field instead.
$ws.Range("D2").Value = '${record.salesPersonFullName}'

# Match the author's final selection/active cell in the sheet.
$ws.Range("D3").Select() | Out-Null
